$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) to reflect the new update.
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell/selection that was saved with the workbook.
$ws.Range("E8").Select()
